$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right after "2021-Q4", cloning its
#    layout/formatting, then writing the new quarter's fund data.
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Copy header/row formatting (styles, borders, column layout) from 2021-Q4.
# (Split the copy so we don't drag along 2021-Q4's untouched/empty A1 cell.)
$q4.Range("B1:H3").Copy($q1.Range("B1"))
$q4.Range("A2:A3").Copy($q1.Range("A2"))

# Data columns must stay text-typed (mirrors the other quarter sheets,
# e.g. "0.73" / "011685" are stored as text, not numbers)
$q1.Range("B2:G3").NumberFormat = "@"

$q1.Range("B2").Value = "011685"
$q1.Range("C2").Value = "创金合信先进装备股票A"
$q1.Range("D2").Value = "0.73"
$q1.Range("E2").Value = "92.01"
$q1.Range("F2").Value = "8.39"
$q1.Range("G2").Value = "0.0612"
$q1.Range("H2").Value = 6

$q1.Range("B3").Value = "011686"
$q1.Range("C3").Value = "创金合信先进装备股票C"
$q1.Range("D3").Value = "0.17"
$q1.Range("E3").Value = "92.01"
$q1.Range("F3").Value = "8.39"
$q1.Range("G3").Value = "0.0143"
$q1.Range("H3").Value = 6

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: push the existing two rows
#    down by one and insert the new 2022-Q1 totals at the top.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift existing data (rows 2:3) down to rows 3:4
$total.Range("A2:D3").Copy($total.Range("A3"))

# New first row: 2022-Q1 totals
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.08

# Re-number the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
